$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.860.33"
$ws.Range("E2").Value = "  +2.03%  "

$ws.Range("D3").Value = "3.717.90"
$ws.Range("E3").Value = "  +1.25%  "

$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "610.59"
$ws.Range("E5").Value = "  +5.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "188.77"
$ws.Range("E6").Value = "  +9.01%  "

$ws.Range("D7").Value = "3.713.76"
$ws.Range("E7").Value = "  +1.43%  "

$ws.Range("E8").Value = "  +1.54%  "

$ws.Range("E9").Value = "  +0.01%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.716"
$ws.Range("E10").Value = "  +1.13%  "

$ws.Range("E11").Value = "  -1.80%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "57.02"
$ws.Range("E12").Value = "  +10.18%  "

$ws.Range("E13").Value = "  -2.40%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.58"
$ws.Range("E14").Value = "  +0.61%  "

$ws.Range("D15").Value = "4.315.20"
$ws.Range("E15").Value = "  +1.23%  "

$ws.Range("D16").Value = "3.720.33"
$ws.Range("E16").Value = "  -0.57%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.31"
$ws.Range("E17").Value = "  +0.26%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.00"
$ws.Range("E18").Value = "  +0.73%  "

$ws.Range("E19").Value = "  -0.57%  "

$ws.Range("E20").Value = "  +0.57%  "

$ws.Range("D21").Value = "68.727.68"
$ws.Range("E21").Value = "  +1.79%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "410.41"
$ws.Range("E22").Value = "  +1.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.58"
$ws.Range("E23").Value = "  -0.39%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "89.14"
$ws.Range("E24").Value = "  +1.31%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.04"
$ws.Range("E25").Value = "  -0.55%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.86"
$ws.Range("E26").Value = "  +1.01%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.03"
$ws.Range("E27").Value = "  +3.82%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.04"
$ws.Range("E28").Value = "  +1.87%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.78"
$ws.Range("E29").Value = "  -0.86%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.64"
$ws.Range("E30").Value = "  +1.64%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.16"
$ws.Range("E31").Value = "  +1.47%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.30"
$ws.Range("E32").Value = "  -8.30%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.123"
$ws.Range("E33").Value = "  +4.67%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.60"
$ws.Range("E34").Value = "  -0.31%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "622.51"
$ws.Range("E35").Value = "  +4.40%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "44.55"
$ws.Range("E36").Value = "  +1.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "65.89"
$ws.Range("E37").Value = "  -0.38%  "

$ws.Range("D38").Value = "0.0₃0833"
$ws.Range("E38").Value = "  -8.96%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.412"
$ws.Range("E39").Value = "  +2.95%  "

$ws.Range("E40").Value = "  +0.22%  "

$ws.Range("E41").Value = "  -0.15%  "

$ws.Range("E42").Value = "  +3.79%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.05"
$ws.Range("E43").Value = "  +0.56%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0443"
$ws.Range("E44").Value = "  +1.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.62"
$ws.Range("E45").Value = "  +2.61%  "

$ws.Range("E46").Value = "  +4.16%  "

$ws.Range("D47").Value = "2.851.49"
$ws.Range("E47").Value = "  +4.23%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.75"
$ws.Range("E48").Value = "  +3.37%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.12"
$ws.Range("E49").Value = "  -3.80%  "

# Row 50/51 coin swap
$ws.Range("B50").Value = "ApeXProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.13"
$ws.Range("E50").Value = "  +0.81%  "

$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.65"
$ws.Range("E51").Value = "  -18.98%  "
